# Atualização da base de faturamento diário (ADD):
# - Corrige valores de total_venda para os dias 21, 24, 25 e 28 de julho/2025
# - Inclui o registro faltante do dia 29/07/2025, deslocando as linhas
#   subsequentes (jun/mai/abr 2025) uma posição para baixo

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige valores existentes (em linha, sem deslocamento)
$ws.Cells.Item(16, 2).Value = 497733.52   # Dia 21/07/2025
$ws.Cells.Item(19, 2).Value = 5000.2      # Dia 24/07/2025
$ws.Cells.Item(20, 2).Value = 7210.32     # Dia 25/07/2025
$ws.Cells.Item(21, 2).Value = 9500.17     # Dia 28/07/2025

# Insere uma nova linha na posição 22 para o dia 29/07/2025,
# empurrando as linhas seguintes para baixo
$ws.Rows.Item(22).Insert()

$ws.Cells.Item(22, 1).Value = 29
$ws.Cells.Item(22, 2).Value = 23549.84
$ws.Cells.Item(22, 3).Value = 7
$ws.Cells.Item(22, 4).Value = 2025
$ws.Cells.Item(22, 5).Value = "07/2025"
